$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6102602746826733
$ws.Range("C2").Value = 0.1727439139589819
$ws.Range("D2").Value = 0.01676340789866515
$ws.Range("F2").Value = 0.4316706243986914
$ws.Range("G2").Value = 0.002377487232050765
$ws.Range("I2").Value = 0.3047803667561588
$ws.Range("M2").Value = 0.8941860217118602
$ws.Range("O2").Value = 1.357150278412661
$ws.Range("B3").Value = 0.5329805854129575
$ws.Range("C3").Value = 0.1542519493905559
$ws.Range("D3").Value = 0.01482873341761604
$ws.Range("F3").Value = 0.4272857750376531
$ws.Range("G3").Value = 0.002380091821955285
$ws.Range("I3").Value = 0.3100628151317508
$ws.Range("M3").Value = 0.7956521778414043
$ws.Range("O3").Value = 1.355696221567882
$ws.Range("B4").Value = 0.4853348552182695
$ws.Range("C4").Value = 0.1428302557792449
$ws.Range("D4").Value = 0.01363448727448002
$ws.Range("F4").Value = 0.424992840530777
$ws.Range("G4").Value = 0.002381775618125186
$ws.Range("I4").Value = 0.3136052725858711
$ws.Range("M4").Value = 0.735678602262297
$ws.Range("O4").Value = 1.356100921375756
$ws.Range("B5").Value = 0.4658714285372128
$ws.Range("C5").Value = 0.1381592904448894
$ws.Range("D5").Value = 0.01314626948086328
$ws.Range("F5").Value = 0.4241585502619287
$ws.Range("G5").Value = 0.002382483111634519
$ws.Range("I5").Value = 0.3151237670265523
$ws.Range("M5").Value = 0.7113651255858997
$ws.Range("O5").Value = 1.356590919026559
$ws.Range("B6").Value = 0.4626367314768345
$ws.Range("C6").Value = 0.13738269509318
$ws.Range("D6").Value = 0.01306510878740141
$ws.Range("F6").Value = 0.4240260531228586
$ws.Range("G6").Value = 0.002382601880828651
$ws.Range("I6").Value = 0.3153804284547341
$ws.Range("M6").Value = 0.7073353377297877
$ws.Range("O6").Value = 1.356691881981035
$ws.Range("B7").Value = 0.4850725541830059
$ws.Range("C7").Value = 0.1427673278430461
$ws.Range("D7").Value = 0.0136279092279068
$ws.Range("F7").Value = 0.4249811841345874
$ws.Range("G7").Value = 0.002381785073031027
$ws.Range("I7").Value = 0.3136254485691303
$ws.Range("M7").Value = 0.7353501991427152
$ws.Range("O7").Value = 1.356106214982901
$ws.Range("B8").Value = 0.5836559076192316
$ws.Range("C8").Value = 0.1663821714001585
$ws.Range("D8").Value = 0.01609767517012983
$ws.Range("F8").Value = 0.4300756366874126
$ws.Range("G8").Value = 0.00237836778357333
$ws.Range("I8").Value = 0.3065395298653328
$ws.Range("M8").Value = 0.8600984758293322
$ws.Range("O8").Value = 1.356378905256804
$ws.Range("B9").Value = 0.7753549698638835
$ws.Range("C9").Value = 0.2121370008084682
$ws.Range("D9").Value = 0.02088882226556876
$ws.Range("F9").Value = 0.4432507971852928
$ws.Range("G9").Value = 0.002372334394318477
$ws.Range("I9").Value = 0.2950288719113381
$ws.Range("M9").Value = 1.109202655633126
$ws.Range("O9").Value = 1.367264631256575
$ws.Range("B10").Value = 0.915126799295706
$ws.Range("C10").Value = 0.2453946683571928
$ws.Range("D10").Value = 0.02437520670458326
$ws.Range("F10").Value = 0.4548962234985083
$ws.Range("G10").Value = 0.002368304482363534
$ws.Range("I10").Value = 0.2880426474980098
$ws.Range("M10").Value = 1.295378651028301
$ws.Range("O10").Value = 1.381653497705685
$ws.Range("B11").Value = 0.9784643512175535
$ws.Range("C11").Value = 0.2604424270064385
$ws.Range("D11").Value = 0.0259535535558868
$ws.Range("F11").Value = 0.4606261142852688
$ws.Range("G11").Value = 0.00236655770939395
$ws.Range("I11").Value = 0.2851873916873657
$ws.Range("M11").Value = 1.380858028415375
$ws.Range("O11").Value = 1.389604757560022
$ws.Range("B12").Value = 1.002411728739844
$ws.Range("C12").Value = 0.2661285110272047
$ws.Range("D12").Value = 0.02655009612938386
$ws.Range("F12").Value = 0.4628584311000452
$ws.Range("G12").Value = 0.002365908613175067
$ws.Range("I12").Value = 0.2841529026104013
$ws.Range("M12").Value = 1.413347815844489
$ws.Range("O12").Value = 1.392819189730147
$ws.Range("B13").Value = 0.9972559161310528
$ws.Range("C13").Value = 0.2649044602527511
$ws.Range("D13").Value = 0.02642167159024211
$ws.Range("F13").Value = 0.4623748748431922
$ws.Range("G13").Value = 0.002366047858565671
$ws.Range("I13").Value = 0.2843736148074321
$ws.Range("M13").Value = 1.406345074486438
$ws.Range("O13").Value = 1.392117834069779
$ws.Range("B14").Value = 0.9804352723217562
$ws.Range("C14").Value = 0.2609104713276906
$ws.Range("D14").Value = 0.02600265464745632
$ws.Range("F14").Value = 0.4608085130281552
$ws.Range("G14").Value = 0.002366504060217704
$ws.Range("I14").Value = 0.285101345220113
$ws.Range("M14").Value = 1.383528522430112
$ws.Range("O14").Value = 1.389865125538137
$ws.Range("B15").Value = 0.9701272444260667
$ws.Range("C15").Value = 0.2584624379080651
$ws.Range("D15").Value = 0.02574584442997718
$ws.Range("F15").Value = 0.4598572254172453
$ws.Range("G15").Value = 0.002366785106099339
$ws.Range("I15").Value = 0.285553196639615
$ws.Range("M15").Value = 1.369568655531467
$ws.Range("O15").Value = 1.388511812719088
$ws.Range("B16").Value = 0.9109824417145092
$ws.Range("C16").Value = 0.244409584317026
$ws.Range("D16").Value = 0.0242719005826828
$ws.Range("F16").Value = 0.4545304922656044
$ws.Range("G16").Value = 0.002368420372410104
$ws.Range("I16").Value = 0.2882357678435881
$ws.Range("M16").Value = 1.289808839894036
$ws.Range("O16").Value = 1.381162263610577
$ws.Range("B17").Value = 0.8746348261060461
$ws.Range("C17").Value = 0.2357674419987177
$ws.Range("D17").Value = 0.02336569726730886
$ws.Range("F17").Value = 0.4513737021488922
$ws.Range("G17").Value = 0.002369445653096449
$ws.Range("I17").Value = 0.2899643418283198
$ws.Range("M17").Value = 1.241085570046152
$ws.Range("O17").Value = 1.377014527866407
$ws.Range("B18").Value = 0.8537056576601003
$ws.Range("C18").Value = 0.2307890846191185
$ws.Range("D18").Value = 0.02284375796472204
$ws.Range("F18").Value = 0.449598671037009
$ws.Range("G18").Value = 0.00237004350871179
$ws.Range("I18").Value = 0.2909889377060715
$ws.Range("M18").Value = 1.213134851996514
$ws.Range("O18").Value = 1.374761091135326
$ws.Range("B19").Value = 0.8466155117849326
$ws.Range("C19").Value = 0.2291022033484751
$ws.Range("D19").Value = 0.02266691687908917
$ws.Range("F19").Value = 0.4490046502188321
$ws.Range("G19").Value = 0.002370247332683009
$ws.Range("I19").Value = 0.2913410540604779
$ws.Range("M19").Value = 1.203683639781147
$ws.Range("O19").Value = 1.374020786609833
$ws.Range("B20").Value = 0.8785064840949985
$ws.Range("C20").Value = 0.2366882055296742
$ws.Range("D20").Value = 0.02346223850631191
$ws.Range("F20").Value = 0.4517055359132485
$ws.Range("G20").Value = 0.00236933566826421
$ws.Range("I20").Value = 0.2897771871545984
$ws.Range("M20").Value = 1.246264572873159
$ws.Range("O20").Value = 1.377442366396252
$ws.Range("B21").Value = 0.9853769272280601
$ws.Range("C21").Value = 0.2620839369618579
$ws.Range("D21").Value = 0.02612576147360102
$ws.Range("F21").Value = 0.461266891480463
$ws.Range("G21").Value = 0.002366369727740809
$ws.Range("I21").Value = 0.2848863222926639
$ws.Range("M21").Value = 1.390226963500126
$ws.Range("O21").Value = 1.390521268151957
$ws.Range("B22").Value = 1.055005361385042
$ws.Range("C22").Value = 0.2786103249115968
$ws.Range("D22").Value = 0.02785984332295044
$ws.Range("F22").Value = 0.4678804401264927
$ws.Range("G22").Value = 0.002364503379829103
$ws.Range("I22").Value = 0.2819624092273116
$ws.Range("M22").Value = 1.485021264385452
$ws.Range("O22").Value = 1.400255638450602
$ws.Range("B23").Value = 1.017863887849671
$ws.Range("C23").Value = 0.2697965531425837
$ws.Range("D23").Value = 0.02693495809737101
$ws.Range("F23").Value = 0.4643171792542518
$ws.Range("G23").Value = 0.002365492911599433
$ws.Range("I23").Value = 0.2834979142835081
$ws.Range("M23").Value = 1.434360603402979
$ws.Range("O23").Value = 1.394951213894558
$ws.Range("B24").Value = 0.8767562084177598
$ws.Range("C24").Value = 0.2362719590591098
$ws.Range("D24").Value = 0.02341859517071754
$ws.Range("F24").Value = 0.4515553897732971
$ws.Range("G24").Value = 0.002369385366243831
$ws.Range("I24").Value = 0.2898617038224636
$ws.Range("M24").Value = 1.243922956444166
$ws.Range("O24").Value = 1.377248532188588
$ws.Range("B25").Value = 0.7236773554402021
$ws.Range("C25").Value = 0.1998205556671451
$ws.Range("D25").Value = 0.01959846411855182
$ws.Range("F25").Value = 0.4393431128729759
$ws.Range("G25").Value = 0.002373895533297976
$ws.Range("I25").Value = 0.2978857639372343
$ws.Range("M25").Value = 1.041290444874917
$ws.Range("O25").Value = 1.363203248816006